$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.075.52"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.24"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.20"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +2.69%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.21"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.73"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "2.661.47"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "2.341.91"
$ws.Range("E16").Value = "  +3.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "42.980.30"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.78"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.75"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.91"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.00"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.24"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.08"
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.81"
$ws.Range("E35").Value = "  +3.52%  "
$ws.Range("E36").Value = "  +2.67%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "2.019.03"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.24"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.57"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("E49").Value = "  -2.46%  "
$ws.Range("D50").Value = "2.529.65"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("E51").Value = "  -1.22%  "
